$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 224, pushing existing rows 224-299 down to 225-300
$ws.Rows(224).Insert()

# Populate the newly inserted row 224 with the new record's data
$ws.Range("A224").Value = 4
$ws.Range("B224").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C224").Value = "Los Lagos"
$ws.Range("D224").Value = 44559
$ws.Range("E224").Value = 10
$ws.Range("F224").Value = 100114001
$ws.Range("G224").Value = "Papa"
$ws.Range("H224").Value = "Pehuenche"
$ws.Range("I224").Value = "1a nueva(o)"
$ws.Range("J224").Value = 400
$ws.Range("K224").Value = 11000
$ws.Range("L224").Value = 12000
$ws.Range("M224").Value = 11500
$ws.Range("N224").Value = "`$/saco 25 kilos"
$ws.Range("O224").Value = "Región de La Araucanía"
$ws.Range("P224").Value = 460
$ws.Range("Q224").Value = 25
$ws.Range("R224").Value = "Hortaliza"
